$d = $word.ActiveDocument
$d.Content.Find.Execute(" Java, Python, microservices, scalable systems", $true, $false, $false, $false, $false, $true, 1, $false, " Java, C#, Python, microservices, scalable and resilient systems", 2)
